# Optuna Attempt (go back with original)
# Updates the forecast comparison values (MyForecast, Inventory Coverage,
# Seasonality Index) on the "Forecast Comparison" sheet, and the derived
# forecast totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Row => (D: MyForecast, H: Inventory Coverage, L: Seasonality Index)
$rows = @{
    2  = @(1, 21, 0.83)
    3  = @(1, 20, 1.09)
    4  = @(1, 19, 0.8100000000000001)
    5  = @(1, 18, 1.01)
    6  = @(1, 17, 1.2)
    7  = @(1, 16, 1.11)
    8  = @(1, 15, 1.11)
    9  = @(1, 14, 1.09)
    10 = @(1, 13, 1.15)
    11 = @(1, 12, 1.04)
    12 = @(1, 11, 1.16)
    13 = @(1, 10, 0.97)
    14 = @(1, 9, 1.15)
    15 = @(1, 8, 1.13)
    16 = @(1, 7, 0.9)
    17 = @(1, 6, 0.98)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $wsForecast.Cells.Item($r, 4).Value = $vals[0]
    $wsForecast.Cells.Item($r, 8).Value = $vals[1]
    $wsForecast.Cells.Item($r, 12).Value = $vals[2]
}

# These look numeric but must stay text cells (matching the existing
# "14" / "7" / "3" text values they replace), so prefix with an
# apostrophe the same way a user would type a text-forced number into
# Excel, then reset the style back to Normal so the quote-prefix flag
# doesn't leave a lingering cell-format change.
$cB9 = $wsSummary.Range("B9")
$cB9.Value = "'16"
$cB9.Style = "Normal"

$cB10 = $wsSummary.Range("B10")
$cB10.Value = "'8"
$cB10.Style = "Normal"

$cB11 = $wsSummary.Range("B11")
$cB11.Value = "'4"
$cB11.Style = "Normal"
